# Runs.xlsx edit — adds the N=30000 / 30sect "Pictures and results" block
# (row 24 header + new annotation cells in rows 25-30) and marks row 3 with
# a "tuloste" note, per the target diff.
#
# NOTE on ordering: new shared-string entries are appended to the workbook's
# shared-string table in first-write order, so the cell writes below are
# deliberately sequenced to reproduce the exact target shared-string index
# order (25: SOA_formation_N_30000 .. 32: ylarajan ylitse).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new shared strings, introduced in the exact order they must land in
#     the shared-string table --------------------------------------------
$ws.Range("A25").Value = "SOA_formation_N_30000"        # sst 25
$ws.Range("B24").Value = "30sect"                        # sst 26
$ws.Range("A26").Value = "run_20130619T183705"           # sst 27
$ws.Range("I3").Value  = "tuloste"                        # sst 28
$ws.Range("K30").Value = "NaN"                            # sst 29
$ws.Range("L30").Value = "jakauma katosi"                 # sst 30
$ws.Range("L25").Value = "Ongelmia jakauman kasvaessa"    # sst 31
$ws.Range("L26").Value = "ylärajan ylitse"                # sst 32

# --- remaining new cells (reuse existing shared strings) -----------------
$ws.Range("A24").Value = "Script"
$ws.Range("I25").Value = "tuloste"
$ws.Range("I27").Value = "tuloste"

$ws.Range("A27").Value = "Pictures"
$ws.Range("A27").Font.Bold = $true

$ws.Range("A28").Value = "deltaP and deltaMoa"
$ws.Range("A29").Value = "Y(t)"
$ws.Range("J29").Value = "no result, matlab tolerance error"
$ws.Range("A30").Value = "Y(deltaMoa)"

# --- E-column toggles: swap the "1/(24*60*60)" text marker for a plain 0
#     and vice versa on alternating rows -----------------------------------
$ws.Range("E25").Value = 0
$ws.Range("E26").Value = "1/(24*60*60)"
$ws.Range("E27").Value = 0
$ws.Range("E28").Value = "1/(24*60*60)"
$ws.Range("E29").Value = 0
$ws.Range("E30").Value = "1/(24*60*60)"

# --- final selection, matching the saved cursor position in the target ---
$ws.Range("L27").Select() | Out-Null
